# Apply irrigation measure changes for the last day (rows 128-145, column F)
# and update the active sheet view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update irrigation values (column F) for rows 128 through 145 from 0 to 50
$ws.Range("F128:F145").Value = 50

# Update the sheet view's top-left visible cell and active selection
$ws.Application.ActiveWindow.ScrollRow = 120
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("J135").Select()
